$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new title row at the very top (pushes the existing "SMT" data,
#    which had no title before, down by one row).
# ---------------------------------------------------------------------------
$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "SMT"
$ws.Rows("1:1").RowHeight = 21
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16
$ws.Range("A1").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 2. Insert 11 fresh rows below the first table (row 19 is now the last row
#    of table 1) to make room for: one blank separator row (20), the new
#    "no SMT" title block (21-24), a blank separator (25) and the second
#    comparison table (26-30).
# ---------------------------------------------------------------------------
$ws.Rows("20:30").Insert()

# --- second title block -----------------------------------------------------
$ws.Range("A21").Value = "no SMT"
$ws.Rows("21:21").RowHeight = 21
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").Font.Size = 16

$ws.Range("A22").Value = "iterations = 100"
$ws.Range("A23").Value = "2 processes/2cores"
$ws.Range("A24").Value = " node index.js >& /dev/null"

# --- second table header (row 26) -------------------------------------------
$ws.Rows("26:26").RowHeight = 47.25

$ws.Range("A26").Style = "Comma"
$ws.Range("F26").Style = "Comma"

$ws.Range("B26").Value = "code no superpage, data no superpage"
$ws.Range("C26").Value = "code no superpage, data superpage"
$ws.Range("D26").Value = "code superpage, data no superpage"
$ws.Range("E26").Value = "code superpage, data superpage"
$ws.Range("G26").Value = "Col C/Col B"
$ws.Range("H26").Value = "Col D/ Col B"
$ws.Range("I26").Value = "Col E/ Col B"
$ws.Range("B26:E26").Style = "Comma"
$ws.Range("G26:I26").Style = "Comma"
$ws.Range("B26:I26").WrapText = $true

# --- second table data rows -------------------------------------------------
$ws.Range("A27").Value = "CPU_CLK_UNHALTED.THREAD_P"
$ws.Range("B27").Value = 18249858809.810001
$ws.Range("C27").Value = 18270850087.305
$ws.Range("D27").Value = 18077693604.705002
$ws.Range("E27").Value = 18046952072.825001
$ws.Range("F27").Style = "Comma"
$ws.Range("G27").Formula = "=C27/B27"
$ws.Range("H27").Formula = "=D27/B27"
$ws.Range("I27").Formula = "=E27/B27"

$ws.Range("A28").Value = "DTLB_LOAD_MISSES.WALK_PENDING"
$ws.Range("B28").Value = 87217879.090000004
$ws.Range("C28").Value = 88646873.629999995
$ws.Range("D28").Value = 82384890.215000004
$ws.Range("E28").Value = 75467389.984999999
$ws.Range("F28").Style = "Comma"
$ws.Range("G28").Formula = "=C28/B28"
$ws.Range("H28").Formula = "=D28/B28"
$ws.Range("I28").Formula = "=E28/B28"

$ws.Range("A29").Value = "DTLB_STORE_MISSES.WALK_PENDING"
$ws.Range("B29").Value = 70270472.165000007
$ws.Range("C29").Value = 70144121.525000006
$ws.Range("D29").Value = 70457187.969999999
$ws.Range("E29").Value = 67742924.480000004
$ws.Range("F29").Style = "Comma"
$ws.Range("G29").Formula = "=C29/B29"
$ws.Range("H29").Formula = "=D29/B29"
$ws.Range("I29").Formula = "=E29/B29"

$ws.Range("A30").Value = "elapse time"
$ws.Range("B30").Value = 1029.3789999999999
$ws.Range("C30").Value = 1029.58
$ws.Range("D30").Value = 1015.258
$ws.Range("E30").Value = 1014.2805
$ws.Range("F30").Style = "Comma"
$ws.Range("G30").Formula = "=C30/B30"
$ws.Range("H30").Formula = "=D30/B30"
$ws.Range("I30").Formula = "=E30/B30"

# Row-label formatting for the second table to mirror the first table:
# plain labels in B8:B29-style font, red for the "elapse time" summary row.
$ws.Range("A27:A29").Font.Name = "Calibri"
$ws.Range("A27:A29").Font.Size = 12
$ws.Range("A30").Font.Color = 255

# ---------------------------------------------------------------------------
# 3. Sheet view: scroll down a bit and move the selection, matching the
#    author's on-save cursor position.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 6
